$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = 1
$ws.Range("C44").Value = "2024-06-15 19:10:42"
$ws.Range("D44").Value = 200
$ws.Range("E44").Value = 16

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 2
$ws.Range("C45").Value = "2024-06-15 19:10:42"
$ws.Range("D45").Value = 200
$ws.Range("E45").Value = 1
